$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column (C) for rows 2-7
# from serial date 45185 (2023-09-16) to 45204 (2023-10-05),
# keeping the existing date number format/style on the cells.
foreach ($row in 2..7) {
    $ws.Cells.Item($row, 3).Value = 45204
}
